$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 'Nom Nom'
$ws.Range("G2").Value = 'Surf Co'
$ws.Range("H2").Value = 'Laddy'
$ws.Range("I2").Value = 'Padge'
$ws.Range("E3").Value = 'Unicorn'
$ws.Range("F3").Value = 'Goobie'
$ws.Range("G3").Value = 'GoGo'
$ws.Range("H3").Value = 'Bow Wow'
$ws.Range("I3").Value = 'Tross'
$ws.Range("E4").Value = 'Hawma'
$ws.Range("I4").Value = 'Surf Co'
$ws.Range("E5").Value = 'Sea'
$ws.Range("I5").Value = 'Nono'
$ws.Range("E6").Value = 'Hullabaloo'
$ws.Range("G6").Value = 'Gaggles'
$ws.Range("H6").Value = 'Gaggles'
$ws.Range("F7").Value = 'GoGo'
$ws.Range("G7").Value = 'Bonez'
$ws.Range("H7").Value = 'G-Poppy'
$ws.Range("I7").Value = 'Sea'
$ws.Range("I9").Value = 'Bow Wow'
$ws.Range("I11").Value = 'Surf Co'
$ws.Range("I12").Value = 'Indi'
$ws.Range("F13").Value = 'Sea'
$ws.Range("E14").Value = 'Surf Co'
$ws.Range("E15").Value = 'Tross'
$ws.Range("F15").Value = 'Tross'
$ws.Range("E16").Value = 'Bonez'
$ws.Range("E17").Value = 'Sea'
$ws.Range("F17").Value = 'Captain'
$ws.Range("H17").Value = 'Tross'
$ws.Range("I17").Value = 'Captain'
$ws.Range("E18").Value = 'Bow Wow'
$ws.Range("F18").Value = 'Blister'
$ws.Range("H18").Value = 'Bow Wow'
$ws.Range("I18").Value = 'Smiles'
$ws.Range("F20").Value = 'Tross'
$ws.Range("F21").Value = 'Blister'
$ws.Range("H22").Value = 'Sea'
$ws.Range("H23").Value = 'Indi'
$ws.Range("E24").Value = 'Indi'
$ws.Range("F24").Value = 'Sea'
$ws.Range("H24").Value = 'Bonez'
$ws.Range("I24").Value = 'Sea'
$ws.Range("F25").Value = 'Bonez'
$ws.Range("G25").Value = 'Surf Co'
$ws.Range("H25").Value = 'Tross'
$ws.Range("I25").Value = 'Indi'
$ws.Range("D26").Value = 'Hoops'
$ws.Range("E26").Value = 'GoGo'
$ws.Range("G26").Value = 'Goobie'
$ws.Range("H26").Value = 'Indi'
$ws.Range("I26").Value = 'Unicorn'
$ws.Range("D27").Value = 'Unicorn'
$ws.Range("E27").Value = 'Nom Nom'
$ws.Range("G27").Value = 'Hoops'
$ws.Range("H27").Value = 'Sea'
$ws.Range("I27").Value = 'Goobie'
$ws.Range("D28").Value = 'Bonez'
$ws.Range("E28").Value = 'Ups'
$ws.Range("F28").Value = 'Indi'
$ws.Range("G28").Value = 'Ups'
$ws.Range("I28").Value = 'Bonez'
$ws.Range("D29").Value = 'Tross'
$ws.Range("E29").Value = 'Tross'
$ws.Range("F29").Value = 'Zombie'
$ws.Range("G29").Value = 'Unicorn'
$ws.Range("H29").Value = 'Ups'
$ws.Range("I29").Value = 'Nom Nom'
$ws.Range("E30").Value = 'Hoops'
$ws.Range("F30").Value = 'Surf Co'
$ws.Range("G30").Value = 'Zombie'
$ws.Range("H30").Value = 'Nom Nom'
$ws.Range("F31").Value = 'Captain'
$ws.Range("H31").Value = 'Padge'
$ws.Range("F32").Value = 'Bow Wow'
$ws.Range("H32").Value = 'Bow Wow'
$ws.Range("E33").Value = 'Chicken'
$ws.Range("F33").Value = 'T-Whisk'
$ws.Range("G33").Value = 'Smiles'
$ws.Range("H33").Value = 'Chicken'
$ws.Range("I33").Value = 'Chicken'
$ws.Range("E34").Value = 'Smiles'
$ws.Range("G34").Value = 'Bow Wow'
$ws.Range("G35").Value = 'Sea'
$ws.Range("H35").Value = 'Smiles'
$ws.Range("I35").Value = 'Bow Wow'
$ws.Range("G36").Value = 'Opps'
$ws.Range("H36").Value = 'Pizza'
$ws.Range("I36").Value = 'Smiles'
$ws.Range("D37").Value = 'Burning Bush'
$ws.Range("E37").Value = 'Nodder'
$ws.Range("F37").Value = 'G-Poppy'
$ws.Range("H37").Value = 'Nono'
$ws.Range("I37").Value = 'Zombie'
$ws.Range("E38").Value = 'Captain'
$ws.Range("G38").Value = 'Padge'
$ws.Range("I38").Value = 'Captain'
$ws.Range("G39").Value = 'Nom Nom'
$ws.Range("I39").Value = 'Hawma'
$ws.Range("E40").Value = 'Laddy'
$ws.Range("G40").Value = 'Stastro'
$ws.Range("I40").Value = 'Stastro'
$ws.Range("E41").Value = 'Blister'
$ws.Range("F41").Value = 'Laddy'
$ws.Range("G41").Value = 'Nono'
$ws.Range("I41").Value = 'Opps'
$ws.Range("E42").Value = 'Hawma'
$ws.Range("F42").Value = 'Burning Bush'
$ws.Range("H42").Value = 'G-Poppy'
$ws.Range("I42").Value = 'G-Poppy'
$ws.Range("E43").Value = 'Burning Bush'
$ws.Range("F43").Value = 'Nono'
$ws.Range("G43").Value = 'Pizza'
$ws.Range("H43").Value = 'Hawma'
$ws.Range("E44").Value = 'NONE FOUND'
$ws.Range("F44").Value = 'Padge'
$ws.Range("H44").Value = 'Burning Bush'
$ws.Range("I44").Value = 'Padge'
$ws.Range("H45").Value = 'Nodder'
$ws.Range("I45").Value = 'Pizza'
$ws.Range("G47").Value = 'GoGo'
$ws.Range("I47").Value = 'Hoops'
$ws.Range("H48").Value = 'Surf Co'
$ws.Range("H49").Value = 'Opps'
$ws.Range("I49").Value = 'GoGo'
$ws.Range("F50").Value = 'NONE FOUND'
$ws.Range("H50").Value = 'Stastro'
$ws.Range("H51").Value = 'Zombie'
$ws.Range("F54").Value = 'GoGo'
$ws.Range("H54").Value = 'Goobie'
$ws.Range("F55").Value = 'Surf Co'
$ws.Range("H55").Value = 'Ups'
$ws.Range("D56").Value = 'Zombie'
$ws.Range("E56").Value = 'Zombie'
$ws.Range("D57").Value = 'Nom Nom'
$ws.Range("E57").Value = 'Ups'
$ws.Range("D58").Value = 'Sea'
$ws.Range("E58").Value = 'Bonez'
$ws.Range("E59").Value = 'Pizza'
$ws.Range("F59").Value = 'Hoops'
$ws.Range("H59").Value = 'Tross'
$ws.Range("D60").Value = 'Bow Wow'
$ws.Range("E60").Value = 'Unicorn'
$ws.Range("F60").Value = 'Sea'
$ws.Range("G60").Value = 'Hawma'
$ws.Range("H60").Value = 'Jaws'
$ws.Range("I60").Value = 'Burning Bush'
$ws.Range("D61").Value = 'Swamp Puppy'
$ws.Range("E61").Value = 'Swamp Puppy'
$ws.Range("F61").Value = 'Socks'
$ws.Range("G61").Value = '"2319"'
$ws.Range("H61").Value = 'Swamp Puppy'
$ws.Range("I61").Value = 'Packs'
$ws.Range("D62").Value = 'Burning Bush'
$ws.Range("E62").Value = 'Burning Bush'
$ws.Range("F62").Value = 'Laddy'
$ws.Range("G62").Value = 'Surf Co'
$ws.Range("H62").Value = 'Surf Co'
$ws.Range("I62").Value = 'Stastro'
